$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "B" column (query) text for the three data rows, replacing the
# recurrence-score filter from "16-20" to "31-35". The content/role of each
# row's query is unchanged; only the filter literal changes.

$caseQuery = $ws.Range("B2").Value2
$caseQuery = $caseQuery -replace '"16-20"', '"31-35"'
$ws.Range("B2").Value2 = $caseQuery

$sampleQuery = $ws.Range("B3").Value2
$sampleQuery = $sampleQuery -replace '"16-20"', '"31-35"'
$ws.Range("B3").Value2 = $sampleQuery

$fileQuery = $ws.Range("B4").Value2
$fileQuery = $fileQuery -replace '"16-20"', '"31-35"'
$ws.Range("B4").Value2 = $fileQuery

# Also update the "C" column (StatQuery), which contains the same filter
# literal, for all three rows.
$statQuery = $ws.Range("C2").Value2
$statQuery = $statQuery -replace '"16-20"', '"31-35"'
$ws.Range("C2").Value2 = $statQuery
$ws.Range("C3").Value2 = $statQuery
$ws.Range("C4").Value2 = $statQuery

# Reflect the user's last selection/view position (B3 active, top-left A3).
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
